# daily auto push: 2026-01-23 18:50 UTC
#
# A new daily-log entry for 2026/01/23 (金, hour 23, value 15) was added
# to the log table on Sheet1. In the original sheet this pushed a new row
# in right before the existing "2026/12/29" block, shifting every row from
# the old row 684 downward by one (old row 725 -> new row 726). Excel's
# InsertRow (which shifts everything below it down, just like Range.Insert
# in the real UI) reproduces that shift for free, so we only need to
# insert one row and fill in its four values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 684; rows 684..725 all shift down to 685..726,
# and the sheet's used range grows from D725 to D726 automatically.
$ws.Rows.Item(684).Insert()

# The date column is stored as plain text ("2026/01/23"), not a real Excel
# date serial. Force text interpretation while writing the value, then
# drop the temporary "@" number format so the cell ends up unstyled, same
# as every other data cell in the column.
$ws.Range("A684").NumberFormat = "@"
$ws.Range("A684").Value = "2026/01/23"
$ws.Range("A684").ClearFormats()

$ws.Range("B684").Value = "金"
$ws.Range("C684").Value = 23
$ws.Range("D684").Value = 15
